# COREESG_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer (A10) from
#    2021-05-13 to 2021-05-14
#  - refresh the Weight / Percent Change figures (columns D & E, rows 2-7)
#
# The sheet ships password-protected (sheetProtection), so we briefly
# unprotect it to write the new values and re-protect it afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- disclaimer text: update the "as of" date -------------------------
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."
$ws.Range("A10").Value = $disclaimer

# --- Weight (D) / Percent Change (E) updates ---------------------------
$ws.Range("D2").Value = 0.2416348192227899
$ws.Range("E2").Value = 0.0192937123169683

$ws.Range("D3").Value = 0.5055597155062046
$ws.Range("E3").Value = 0.009741969457609523

$ws.Range("D4").Value = 0.09383283268451297
$ws.Range("E4").Value = 0.02688834154351394

$ws.Range("D5").Value = 0.1029536910968928
$ws.Range("E5").Value = 0.01444141689373257

$ws.Range("D6").Value = 0.05601894148959962
$ws.Range("E6").Value = 0.02201622247972179

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = 0.01483031190131512

# restore sheet protection
$ws.Protect()
